$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "SK증권제11호스팩" row (row 15) entirely; the row below it
# (비엔케이제2호스팩) shifts up to take its place.
$ws.Rows.Item(15).Delete()
